$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update the "Total" row correct/total marks figures
$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 65
$ws.Range("E12").Value = "65/140"
